# Added UNIQUE constraints to db definition, changed classification_definition.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting that's already used for the other "constraint" cells
# in this column (e.g. C3, which uses the same Arial/10pt, no border style)
# onto E3, then set its value to the newly introduced "UNIQUE" constraint.
$ws.Range("C3").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("E3").Value = "UNIQUE"

# The active cell/selection moved from D9 (off the used range) to E4.
$ws.Range("E4").Select()
